$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") values on rows 2-92 were recorded as 40
# but should be 240. Update each cell accordingly.
$ws.Range("B2:B92").Value = 240
